$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 79.61246300000001
$ws.Range("H2").Value = 238.837389
$ws.Range("I2").Value = 0.6728436998494041
$ws.Range("J2").Value = 0.6728436998494042
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 98.946724
$ws.Range("N2").Value = 296.840172
$ws.Range("O2").Value = 0.2098009692989996
$ws.Range("P2").Value = 0.2098009692989996
$ws.Range("Q2").Value = 7877.392403421213
$ws.Range("R2").Value = 70896.53163079091
$ws.Range("S2").Value = 0.1411632604151301
$ws.Range("T2").Value = 0.1411632604151302

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 79.61246300000001
$ws.Range("H3").Value = 238.837389
$ws.Range("I3").Value = 0.6728436998494041
$ws.Range("J3").Value = 0.6728436998494042
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 163.0062356666667
$ws.Range("N3").Value = 489.018707
$ws.Range("O3").Value = 0.345629090707923
$ws.Range("P3").Value = 0.3456290907079231
$ws.Range("Q3").Value = 12977.32790578178
$ws.Range("R3").Value = 116795.951152036
$ws.Range("S3").Value = 0.2325543561675042
$ws.Range("T3").Value = 0.2325543561675043

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 79.61246300000001
$ws.Range("H4").Value = 238.837389
$ws.Range("I4").Value = 0.6728436998494041
$ws.Range("J4").Value = 0.6728436998494042
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 65.39610666666668
$ws.Range("N4").Value = 196.18832
$ws.Range("O4").Value = 0.1386621609326595
$ws.Range("P4").Value = 0.1386621609326595
$ws.Range("Q4").Value = 5206.345122344055
$ws.Range("R4").Value = 46857.10610109649
$ws.Range("S4").Value = 0.09329796139104408
$ws.Range("T4").Value = 0.09329796139104411

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 79.61246300000001
$ws.Range("H5").Value = 238.837389
$ws.Range("I5").Value = 0.6728436998494041
$ws.Range("J5").Value = 0.6728436998494042
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 144.2727966666667
$ws.Range("N5").Value = 432.81839
$ws.Range("O5").Value = 0.3059077790604178
$ws.Range("P5").Value = 0.3059077790604179
$ws.Range("Q5").Value = 11485.91268653152
$ws.Range("R5").Value = 103373.2141787837
$ws.Range("S5").Value = 0.2058281218757256
$ws.Range("T5").Value = 0.2058281218757257

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 3.815058666666667
$ws.Range("H6").Value = 11.445176
$ws.Range("I6").Value = 0.03224291890608301
$ws.Range("J6").Value = 0.03224291890608302
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 98.946724
$ws.Range("N6").Value = 296.840172
$ws.Range("O6").Value = 0.2098009692989996
$ws.Range("P6").Value = 0.2098009692989996
$ws.Range("Q6").Value = 377.4875569344747
$ws.Range("R6").Value = 3397.388012410272
$ws.Range("S6").Value = 0.006764595639525255
$ws.Range("T6").Value = 0.006764595639525257

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 3.815058666666667
$ws.Range("H7").Value = 11.445176
$ws.Range("I7").Value = 0.03224291890608301
$ws.Range("J7").Value = 0.03224291890608302
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 163.0062356666667
$ws.Range("N7").Value = 489.018707
$ws.Range("O7").Value = 0.345629090707923
$ws.Range("P7").Value = 0.3456290907079231
$ws.Range("Q7").Value = 621.8783521008259
$ws.Range("R7").Value = 5596.905168907432
$ws.Range("S7").Value = 0.01114409074327877
$ws.Range("T7").Value = 0.01114409074327877

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3.815058666666667
$ws.Range("H8").Value = 11.445176
$ws.Range("I8").Value = 0.03224291890608301
$ws.Range("J8").Value = 0.03224291890608302
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 65.39610666666668
$ws.Range("N8").Value = 196.18832
$ws.Range("O8").Value = 0.1386621609326595
$ws.Range("P8").Value = 0.1386621609326595
$ws.Range("Q8").Value = 249.4899835049245
$ws.Range("R8").Value = 2245.40985154432
$ws.Range("S8").Value = 0.004470872810293971
$ws.Range("T8").Value = 0.004470872810293972

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3.815058666666667
$ws.Range("H9").Value = 11.445176
$ws.Range("I9").Value = 0.03224291890608301
$ws.Range("J9").Value = 0.03224291890608302
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 144.2727966666667
$ws.Range("N9").Value = 432.81839
$ws.Range("O9").Value = 0.3059077790604178
$ws.Range("P9").Value = 0.3059077790604179
$ws.Range("Q9").Value = 550.4091832874044
$ws.Range("R9").Value = 4953.68264958664
$ws.Range("S9").Value = 0.00986335971298501
$ws.Range("T9").Value = 0.009863359712985013

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 30.45313966666667
$ws.Range("H10").Value = 91.359419
$ws.Range("I10").Value = 0.2573743154429307
$ws.Range("J10").Value = 0.2573743154429307
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 98.946724
$ws.Range("N10").Value = 296.840172
$ws.Range("O10").Value = 0.2098009692989996
$ws.Range("P10").Value = 0.2098009692989996
$ws.Range("Q10").Value = 3013.238405531119
$ws.Range("R10").Value = 27119.14564978007
$ws.Range("S10").Value = 0.05399738085259333
$ws.Range("T10").Value = 0.05399738085259335

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 30.45313966666667
$ws.Range("H11").Value = 91.359419
$ws.Range("I11").Value = 0.2573743154429307
$ws.Range("J11").Value = 0.2573743154429307
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 163.0062356666667
$ws.Range("N11").Value = 489.018707
$ws.Range("O11").Value = 0.345629090707923
$ws.Range("P11").Value = 0.3456290907079231
$ws.Range("Q11").Value = 4964.051661294582
$ws.Range("R11").Value = 44676.46495165124
$ws.Range("S11").Value = 0.08895605061811428
$ws.Range("T11").Value = 0.0889560506181143

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 30.45313966666667
$ws.Range("H12").Value = 91.359419
$ws.Range("I12").Value = 0.2573743154429307
$ws.Range("J12").Value = 0.2573743154429307
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 65.39610666666668
$ws.Range("N12").Value = 196.18832
$ws.Range("O12").Value = 0.1386621609326595
$ws.Range("P12").Value = 0.1386621609326595
$ws.Range("Q12").Value = 1991.516769976232
$ws.Range("R12").Value = 17923.65092978608
$ws.Range("S12").Value = 0.03568807874788071
$ws.Range("T12").Value = 0.03568807874788073

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 30.45313966666667
$ws.Range("H13").Value = 91.359419
$ws.Range("I13").Value = 0.2573743154429307
$ws.Range("J13").Value = 0.2573743154429307
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 144.2727966666667
$ws.Range("N13").Value = 432.81839
$ws.Range("O13").Value = 0.3059077790604178
$ws.Range("P13").Value = 0.3059077790604179
$ws.Range("Q13").Value = 4393.559626990601
$ws.Range("R13").Value = 39542.03664291542
$ws.Range("S13").Value = 0.07873280522434231
$ws.Range("T13").Value = 0.07873280522434235

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 4.441711333333334
$ws.Range("H14").Value = 13.325134
$ws.Range("I14").Value = 0.03753906580158222
$ws.Range("J14").Value = 0.03753906580158223
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 98.946724
$ws.Range("N14").Value = 296.840172
$ws.Range("O14").Value = 0.2098009692989996
$ws.Range("P14").Value = 0.2098009692989996
$ws.Range("Q14").Value = 439.4927853870054
$ws.Range("R14").Value = 3955.435068483048
$ws.Range("S14").Value = 0.007875732391750876
$ws.Range("T14").Value = 0.007875732391750878

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 4.441711333333334
$ws.Range("H15").Value = 13.325134
$ws.Range("I15").Value = 0.03753906580158222
$ws.Range("J15").Value = 0.03753906580158223
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 163.0062356666667
$ws.Range("N15").Value = 489.018707
$ws.Range("O15").Value = 0.345629090707923
$ws.Range("P15").Value = 0.3456290907079231
$ws.Range("Q15").Value = 724.0266443646377
$ws.Range("R15").Value = 6516.239799281739
$ws.Range("S15").Value = 0.01297459317902575
$ws.Range("T15").Value = 0.01297459317902576

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 4.441711333333334
$ws.Range("H16").Value = 13.325134
$ws.Range("I16").Value = 0.03753906580158222
$ws.Range("J16").Value = 0.03753906580158223
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 65.39610666666668
$ws.Range("N16").Value = 196.18832
$ws.Range("O16").Value = 0.1386621609326595
$ws.Range("P16").Value = 0.1386621609326595
$ws.Range("Q16").Value = 290.470628137209
$ws.Range("R16").Value = 2614.235653234881
$ws.Range("S16").Value = 0.005205247983440687
$ws.Range("T16").Value = 0.005205247983440689

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 4.441711333333334
$ws.Range("H17").Value = 13.325134
$ws.Range("I17").Value = 0.03753906580158222
$ws.Range("J17").Value = 0.03753906580158223
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 144.2727966666667
$ws.Range("N17").Value = 432.81839
$ws.Range("O17").Value = 0.3059077790604178
$ws.Range("P17").Value = 0.3059077790604179
$ws.Range("Q17").Value = 640.818116046029
$ws.Range("R17").Value = 5767.363044414261
$ws.Range("S17").Value = 0.0114834922473649
$ws.Range("T17").Value = 0.0114834922473649
